$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert two new rows before the current row 200 (shifts rows 200+ down by 2)
$ws.Rows("200:201").Insert()

# Populate the two new rows with the new pin data (order matches the
# shared-string table layout: P27, P28, BL_CLK, PIN_ARRAY_1x1)
$ws.Range("A200").Value = "P27"
$ws.Range("A201").Value = "P28"
$ws.Range("B200").Value = "BL_CLK"
$ws.Range("C200").Value = "PIN_ARRAY_1x1"
$ws.Range("B201").Value = "RST_N"
$ws.Range("C201").Value = "PIN_ARRAY_1x1"

# Update the named range to reflect the two additional rows
$wb.Names.Item("XMP16_03").RefersTo = "=Sheet1!`$A`$9:`$C`$397"

# Keep the selection in sync with where the user was working
$ws.Range("C201").Select()
